$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.376.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.619.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.05"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.38"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.211"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.647"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.26"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000304"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.56"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.198.38"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "604.22"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.477.93"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.613.67"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.09"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.14"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.17"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.61"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.99"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.65"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.79"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.71"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +7.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.29"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.27"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.46"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0881"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.904.65"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.89%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "518.13"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.84"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.53"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0459"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.51"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.55"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.52%  "
